$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginCredentials")

$ws.Range("A3").Value = "minal"
$ws.Range("B3").Value = "hghgbvc"

$ws.Range("B3").Select()
